$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the Qualification values (D2:D9) that previously held 34234
$ws.Range("D2:D9").ClearContents()

# E8 ("Project Number" column, row 8) now references the new quoted value
$ws.Range("E8").Value = '"3211/3441"'

# Move the active selection from A9 to D1
$ws.Range("D1").Select()
